# Lab 5 - Completed upto Q4 Project - Handled services not in list
#
# This script reproduces the authoring session: a new "service not found"
# fallback row (ERROR) was added at the bottom of the table, and a new
# general "G000" info row was inserted right under the header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Add the fallback "service not in list" row at the bottom (row 11),
#    keeping the existing rows 2-9 (G001..G008) where they are for now.
$ws.Range("A11").Value = "ERROR"
$ws.Range("B11").Value = "Sorry I'm not aware about this service. Hence I will inform this to the management and will reach back to you."

# 2) Shift the existing G001..G008 rows down by one (rows 2-9 -> 3-10) to
#    make room for the new general "G000" info row, then insert it.
for ($r = 9; $r -ge 2; $r--) {
    $srcA = $ws.Range("A$r").Value()
    $srcB = $ws.Range("B$r").Value()
    $ws.Range("A" + ($r + 1)).Value = $srcA
    $ws.Range("B" + ($r + 1)).Value = $srcB
}

$ws.Range("A2").Value = "G000"
$ws.Range("B2").Value = "This is Info about saloon services in general!"

# 3) Widen column B so the long descriptions are easier to read.
$ws.Columns.Item(2).ColumnWidth = 151.76

# 4) Leave the cursor where the author left it when saving.
$ws.Range("B15").Select()
